# fdo#51601 test fixture update:
#  - add Sheet2 (after Sheet1) containing a formula that evaluates to an error
#  - add "Fdo#51601" text on Sheet1!B3
#  - move the selection on Sheet1 to D4
#  - mark Sheet1's (and the new Sheet2's) rows as having an explicit
#    (custom) row height, matching the resaved fixture
#  - give the first style font an explicit charset

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1 edits -------------------------------------------------------

# New cell B3 with the bug-tracker reference text
$ws1.Range("B3").Value = "Fdo#51601"

# Touch every used row's height so it round-trips with an explicit
# (custom) row-height flag, without actually changing the visible height.
for ($i = 1; $i -le 5; $i++) {
    $row = $ws1.Rows.Item($i)
    $row.RowHeight = $row.RowHeight
}

# Give the workbook's base font an explicit charset (ANSI/Default = 1)
$ws1.Range("A1").Font.Charset = 1

# --- Add Sheet2 ----------------------------------------------------------

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# A formula referencing a single range (wrong argument count for
# AVERAGEIF) so the cell evaluates to an error value.
$ws2.Range("C1").Formula = "=averageif(A2:B2)"

$row2 = $ws2.Rows.Item(1)
$row2.RowHeight = $row2.RowHeight

[void]$ws2.Range("C1").Select()

# Restore Sheet1 as the active sheet, with the selection on D4
[void]$ws1.Activate()
[void]$ws1.Range("D4").Select()
